$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'298.79"
$ws.Range("E2").Value = "'0.63%"
$ws.Range("D3").Value = "'31.14"
$ws.Range("D4").Value = "'5.109"
$ws.Range("E4").Value = "'-0.21%"
$ws.Range("D5").Value = "'0.08034"
$ws.Range("E5").Value = "'9.67%"
$ws.Range("D6").Value = "'2.593"
$ws.Range("E6").Value = "'50.81%"
$ws.Range("D7").Value = "'7.807"
$ws.Range("E7").Value = "'0.84%"
$ws.Range("D8").Value = "'3.823"
$ws.Range("E8").Value = "'2.72%"
$ws.Range("D9").Value = "'0.9210"
$ws.Range("E9").Value = "'-0.19%"
$ws.Range("D10").Value = "'0.1732"
$ws.Range("E10").Value = "'3.76%"
$ws.Range("D11").Value = "'0.07341"
$ws.Range("E11").Value = "'6.05%"
$ws.Range("D12").Value = "'0.08692"
$ws.Range("E12").Value = "'8.22%"
$ws.Range("D13").Value = "'0.03034"
$ws.Range("E13").Value = "'1.61%"
$ws.Range("D14").Value = "'0.09974"
$ws.Range("E14").Value = "'0.61%"
$ws.Range("D15").Value = "'0.001505"
$ws.Range("E15").Value = "'-0.23%"
$ws.Range("D16").Value = "'0.006100"
$ws.Range("E16").Value = "'-0.02%"
$ws.Range("D17").Value = "'3.507"
$ws.Range("E17").Value = "'1.52%"
$ws.Range("D18").Value = "'2.254"
$ws.Range("E18").Value = "'1.51%"
$ws.Range("D19").Value = "'0.3285"
$ws.Range("E19").Value = "'0.40%"
$ws.Range("D20").Value = "'0.1339"
$ws.Range("E20").Value = "'0.47%"
$ws.Range("D21").Value = "'4.584"
$ws.Range("E21").Value = "'0.59%"
$ws.Range("D22").Value = "'0.1617"
$ws.Range("E22").Value = "'2.29%"
$ws.Range("D23").Value = "'0.04614"
$ws.Range("E23").Value = "'-0.72%"
$ws.Range("D24").Value = "'0.001243"
$ws.Range("E24").Value = "'1.85%"
$ws.Range("D25").Value = "'0.004439"
$ws.Range("E25").Value = "'-6.42%"
$ws.Range("D26").Value = "'0.0001200"
$ws.Range("E26").Value = "'-7.63%"
$ws.Range("D27").Value = "'0.0003432"
$ws.Range("E27").Value = "'83.30%"
$ws.Range("D39").Value = "'0.01813"
$ws.Range("E39").Value = "'6.58%"
$ws.Range("D40").Value = "'0.04530"
$ws.Range("E40").Value = "'2.00%"
$ws.Range("D41").Value = "'0.007085"
$ws.Range("E41").Value = "'-1.57%"
$ws.Range("D42").Value = "'0.1341"
$ws.Range("E42").Value = "'0.91%"
$ws.Range("E43").Value = "'2.34%"
$ws.Range("D44").Value = "'0.009844"
$ws.Range("E44").Value = "'-7.76%"
$ws.Range("D45").Value = "'0.00006747"
$ws.Range("E45").Value = "'11.74%"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("E46").Value = "'0.13%"
$ws.Range("E47").Value = "'-55.60%"
$ws.Range("D49").Value = "'0.00002103"
$ws.Range("E49").Value = "'0.13%"
$ws.Range("D50").Value = "'0.0002003"
$ws.Range("E50").Value = "'0.20%"

Write-Output "Updated 71 cells"
